{"js": "// Replace the date line and each \"N\u00d7D=\" multiplication prompt in the table\n// with the new values from the commit, using a literal find-and-replace\n// over the whole document body (each old value is unique in the doc).\nconst replacements = [\n  [\"2025-01-01 Wednesday\", \"2025-01-02 Thursday\"],\n  [\"522\u00d79=\", \"284\u00d79=\"],\n  [\"182\u00d79=\", \"771\u00d78=\"],\n  [\"251\u00d73=\", \"694\u00d73=\"],\n  [\"912\u00d76=\", \"339\u00d76=\"],\n  [\"363\u00d74=\", \"590\u00d79=\"],\n  [\"988\u00d78=\", \"700\u00d76=\"],\n  [\"613\u00d72=\", \"794\u00d74=\"],\n  [\"312\u00d72=\", \"495\u00d79=\"],\n  [\"257\u00d77=\", \"453\u00d76=\"],\n  [\"990\u00d72=\", \"710\u00d79=\"],\n  [\"291\u00d78=\", \"677\u00d74=\"],\n  [\"842\u00d79=\", \"219\u00d75=\"],\n  [\"826\u00d72=\", \"284\u00d72=\"],\n  [\"930\u00d75=\", \"846\u00d73=\"],\n  [\"395\u00d77=\", \"404\u00d72=\"],\n  [\"352\u00d78=\", \"638\u00d77=\"],\n  [\"291\u00d76=\", \"768\u00d78=\"],\n  [\"702\u00d77=\", \"583\u00d75=\"],\n  [\"657\u00d72=\", \"167\u00d79=\"],\n  [\"601\u00d77=\", \"153\u00d72=\"],\n  [\"568\u00d78=\", \"228\u00d75=\"],\n  [\"859\u00d76=\", \"583\u00d78=\"],\n  [\"835\u00d79=\", \"674\u00d73=\"],\n  [\"323\u00d73=\", \"758\u00d79=\"],\n  [\"735\u00d77=\", \"522\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each \"N\u00d7D=\" multiplication prompt in the table\n# with the new values from the commit, using Find/Replace over the whole\n# document body (each old value is unique in the doc).\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-01-01 Wednesday\", \"2025-01-02 Thursday\"),\n  @(\"522\u00d79=\", \"284\u00d79=\"),\n  @(\"182\u00d79=\", \"771\u00d78=\"),\n  @(\"251\u00d73=\", \"694\u00d73=\"),\n  @(\"912\u00d76=\", \"339\u00d76=\"),\n  @(\"363\u00d74=\", \"590\u00d79=\"),\n  @(\"988\u00d78=\", \"700\u00d76=\"),\n  @(\"613\u00d72=\", \"794\u00d74=\"),\n  @(\"312\u00d72=\", \"495\u00d79=\"),\n  @(\"257\u00d77=\", \"453\u00d76=\"),\n  @(\"990\u00d72=\", \"710\u00d79=\"),\n  @(\"291\u00d78=\", \"677\u00d74=\"),\n  @(\"842\u00d79=\", \"219\u00d75=\"),\n  @(\"826\u00d72=\", \"284\u00d72=\"),\n  @(\"930\u00d75=\", \"846\u00d73=\"),\n  @(\"395\u00d77=\", \"404\u00d72=\"),\n  @(\"352\u00d78=\", \"638\u00d77=\"),\n  @(\"291\u00d76=\", \"768\u00d78=\"),\n  @(\"702\u00d77=\", \"583\u00d75=\"),\n  @(\"657\u00d72=\", \"167\u00d79=\"),\n  @(\"601\u00d77=\", \"153\u00d72=\"),\n  @(\"568\u00d78=\", \"228\u00d75=\"),\n  @(\"859\u00d76=\", \"583\u00d78=\"),\n  @(\"835\u00d79=\", \"674\u00d73=\"),\n  @(\"323\u00d73=\", \"758\u00d79=\"),\n  @(\"735\u00d77=\", \"522\u00d73=\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
